$wb = $excel.ActiveWorkbook

# --- JourneyScreenTwo (sheet4): move the selection from E3 to I2 ---
$ws4 = $wb.Worksheets.Item("JourneyScreenTwo")
$ws4.Activate() | Out-Null
$ws4.Range("I2").Select() | Out-Null

# --- Add the new "JourneyScreenThree" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "JourneyScreenThree"

$ws6.Range("A1").Value = "Father Name"
$ws6.Range("B1").Value = "Mother Name"
$ws6.Range("C1").Value = "Income"
$ws6.Range("A2").Value = "Mr. N. K. Gupta"
$ws6.Range("B2").Value = "Mrs. M K. Gupta"

# Income must stay a text value ("100000000"), not be auto-converted to a number.
$ws6.Range("C2").NumberFormat = "@"
$ws6.Range("C2").Value = "100000000"
$ws6.Range("C2").NumberFormat = "General"

$ws6.Range("C3").Select() | Out-Null
$ws6.Activate() | Out-Null
